$d = $word.ActiveDocument

$replacements = @(
    @("2025-10-15 Wednesday", "2025-10-16 Thursday"),
    @("387÷3=", "571÷2="),
    @("829÷9=", "385÷6="),
    @("729÷7=", "544÷5="),
    @("563÷4=", "516÷9="),
    @("906÷6=", "702÷8="),
    @("359÷6=", "574÷4="),
    @("293÷4=", "350÷8="),
    @("237÷2=", "993÷2="),
    @("415÷6=", "900÷4="),
    @("433÷7=", "578÷8="),
    @("503÷6=", "429÷2="),
    @("335÷3=", "874÷2="),
    @("290÷2=", "142÷8="),
    @("533÷9=", "310÷6="),
    @("741÷3=", "646÷6="),
    @("600÷8=", "509÷4="),
    @("267÷8=", "196÷4="),
    @("755÷8=", "555÷2="),
    @("613÷5=", "471÷8="),
    @("530÷3=", "542÷8="),
    @("521÷5=", "648÷2="),
    @("118÷5=", "890÷3="),
    @("856÷2=", "128÷6="),
    @("449÷2=", "732÷7="),
    @("113÷9=", "588÷7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
